# Edit the "차량" (vehicles) template worksheet:
#  - widen/narrow a few columns
#  - turn row 2 into a concrete example row
#  - append 3 blank data-entry rows (3-5) pre-filled with default numeric values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes -------------------------------------------------
# Excel's ColumnWidth property is offset from the raw OOXML <col width="..">
# value by the default column padding (~0.8333 chars), so subtract that
# offset to land exactly on the target widths from the diff.
$padding = 0.8333333333333334

$ws.Columns.Item(1).ColumnWidth  = 9  - $padding   # A: 6  -> 9
$ws.Columns.Item(2).ColumnWidth  = 9  - $padding   # B: 6  -> 9
$ws.Columns.Item(3).ColumnWidth  = 6  - $padding   # C: 22 -> 6
$ws.Columns.Item(4).ColumnWidth  = 16 - $padding   # D: 11 -> 16
$ws.Columns.Item(17).ColumnWidth = 11 - $padding   # Q: 7  -> 11

# --- Row 2: turn the placeholder row into a filled-in example row ---------
$ws.Range("A2").Value = "예시-V001"
$ws.Range("B2").Value = "12가3456"
$ws.Range("C2").Value = "냉동"
$ws.Range("D2").Value = "UVIS-DVC-12345"
$ws.Range("Q2").Value = "서울특별시 강서구"

# --- Rows 3-5: new blank data-entry rows with default sample values -------
for ($r = 3; $r -le 5; $r++) {
    $ws.Cells.Item($r, 5).Value  = 20      # E: 최대팔레트
    $ws.Cells.Item($r, 6).Value  = 5000    # F: 최대중량(kg)
    $ws.Cells.Item($r, 7).Value  = 30      # G: 최대용적(CBM)
    $ws.Cells.Item($r, 8).Value  = 5       # H: 톤수
    $ws.Cells.Item($r, 9).Value  = 6       # I: 적재함길이(m)
    $ws.Cells.Item($r, 10).Value = 2.4     # J: 적재함너비(m)
    $ws.Cells.Item($r, 11).Value = 2.5     # K: 적재함높이(m)
    $ws.Cells.Item($r, 12).Value = -25     # L: 최저온도
    $ws.Cells.Item($r, 13).Value = -18     # M: 최고온도
    $ws.Cells.Item($r, 14).Value = 5       # N: 연비(km/L)
    $ws.Cells.Item($r, 15).Value = 1500    # O: 리터당연료비
    $ws.Cells.Item($r, 16).Value = "운행가능"  # P: 차량상태
}
